# Add a new hyperlinked text box (a repo link) to slide 1, mirroring the
# author's "Add files via upload" commit which inserted a TextBox shape
# (id=5, name="TextBox 4") right after the existing two placeholder shapes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# PowerPoint's internal shape-id/name counter keeps incrementing even across
# deletes, which is how the original deck ended up with id="5"/"TextBox 4"
# instead of id="4"/"TextBox 3" for the first textbox added to this slide.
# Create-then-delete a throwaway textbox first so the real one picks up the
# expected id/name.
$placeholderShape = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$placeholderShape.Delete()

# Position/size in points (EMU / 12700) matching the target <a:off>/<a:ext>.
$left = 3407943 / 12700
$top = 5273013 / 12700
$width = 8443161 / 12700
$height = 369332 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1

$url = "https://github.com/dougsuh/artCoding/tree/master/BrickBreak"
$tb.TextFrame.TextRange.Text = $url + " "

$urlRange = $tb.TextFrame.TextRange.Characters(1, $url.Length)
$urlRange.ActionSettings(1).Hyperlink.Address = $url
